$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list with refreshed prices / 1h volume percentages.
# For column D (Price) we temporarily force a Text number format before
# assigning the value, then restore the default "Normal" style. This
# prevents Excel from auto-converting numeric-looking strings (e.g.
# "1.00" or "0.0000167") into actual numbers, which would lose the
# original text formatting/precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.931.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.382"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.638.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000167"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.872.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.109.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.508"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0931"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("E37").Value = "  +6.21%  "
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.147.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.666"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("E44").Value = "  +5.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.278.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0257"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.960"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.753"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "263.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.15%  "
